$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 381; this shifts the existing rows 381-436
# down to 382-437, preserving all of their data and formatting.
$ws.Rows.Item(381).Insert()

# Populate the newly inserted row 381 with its data.
$ws.Cells.Item(381, 1).Value = 4
$ws.Cells.Item(381, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(381, 3).Value = "Los Lagos"
$ws.Cells.Item(381, 4).Value = 44984
$ws.Cells.Item(381, 5).Value = 10
$ws.Cells.Item(381, 6).Value = 100112045
$ws.Cells.Item(381, 7).Value = "Zapallo"
$ws.Cells.Item(381, 8).Value = "Paine"
$ws.Cells.Item(381, 9).Value = "1a (cosecha)"
$ws.Cells.Item(381, 10).Value = 250
$ws.Cells.Item(381, 11).Value = 500
$ws.Cells.Item(381, 12).Value = 600
$ws.Cells.Item(381, 13).Value = 540
$ws.Cells.Item(381, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(381, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(381, 16).Value = 540
$ws.Cells.Item(381, 17).Value = 1
$ws.Cells.Item(381, 18).Value = "Hortaliza"
